# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" on all sheets,
# and refresh the related handoff timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2, F2) and the
# "Latest HO Xliff Generate Date" column (G2).
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-10-19 15:21:29"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-10-19 15:21:06"

# de-de sheet: Status (C2) and Latest Handoff Datetime (H2)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-10-19 15:21:29"

# Widen the status columns now that the text is longer (matches Excel's
# own recompute of the "best fit" width for the new text).
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
